$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "02/11/2026"
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 9231.27
$ws.Cells.Item($row, 3).Value = 0.237570868307493
$ws.Cells.Item($row, 4).Value = 0.762429131692507
$ws.Cells.Item($row, 5).Value = -321.18
$ws.Cells.Item($row, 6).Value = -38.46
$ws.Cells.Item($row, 7).Value = -23853.57
$ws.Cells.Item($row, 8).Value = -77.22
$ws.Cells.Item($row, 9).Value = -1159.77
$ws.Cells.Item($row, 10).Value = -34.59
$ws.Cells.Item($row, 11).Value = -25013.34
$ws.Cells.Item($row, 12).Value = -73.04000000000001
